$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) date serial value from 45175 to 45183 for rows 2-8
$ws.Range("C2:C8").Value = 45183
